# 自动更新Excel文件 - 每日巡检脚本
# 规则：
#   E 列(剩余天数) 每天自动减 1；
#   当剩余天数减到 0（即原值为 1）时，视为到期续期：
#       剩余天数重置为 10，开始时间(F列, 格式 yyyymmdd) 顺延 10 天。
#   若 F 列不是合法的 8 位 yyyymmdd 日期（例如数据异常行），则跳过该行不做任何改动。

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)   # E 列：剩余
    $fCell = $ws.Cells.Item($r, 6)   # F 列：开始时间 (yyyymmdd)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = ([int64]$fVal).ToString()
    if ($fStr.Length -ne 8) {
        # 日期格式异常，跳过该行
        continue
    }

    $year = [int]$fStr.Substring(0, 4)
    $month = [int]$fStr.Substring(4, 2)
    $day = [int]$fStr.Substring(6, 2)

    if ($month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    $startDate = Get-Date -Year $year -Month $month -Day $day

    $eNum = [int]$eVal

    if ($eNum -le 1) {
        # 到期，续期 10 天
        $newStart = $startDate.AddDays(10)
        $newFVal = [int]$newStart.ToString("yyyyMMdd")
        $fCell.Value = $newFVal
        $eCell.Value = 10
    } else {
        $eCell.Value = $eNum - 1
    }
}
